# --- add 2022-Q4 data ---------------------------------------------------
# Current layout: [1]="总计"  [2]="2022-Q3" (quarterly fund-holding detail)
#
# Target layout:  [1]="总计"  [2]="2022-Q4" (new detail)  [3]="2022-Q3" (old detail, unchanged)

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ3    = $wb.Worksheets.Item(2)

# 1) Clone the existing "2022-Q3" detail sheet so the old data/formatting is
#    preserved verbatim in its own sheet, placed right after the original.
$wsQ3.Copy($null, $wsQ3)

$wsQ4  = $wb.Worksheets.Item(2)   # the original sheet -> becomes "2022-Q4"
$wsQ3b = $wb.Worksheets.Item(3)   # the clone -> stays "2022-Q3"

$wsQ4.Name  = "2022-Q4"
$wsQ3b.Name = "2022-Q3"

# 2) Re-style the "2022-Q4" sheet to match the "总计" header look & page
#    margins, then overwrite its data with the Q4 numbers.
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats

$wsQ4.PageSetup.LeftMargin   = 54
$wsQ4.PageSetup.RightMargin  = 54
$wsQ4.PageSetup.TopMargin    = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

function Set-TextCell($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $text
    $c.ClearFormats()
}

# Row 2 (fund 008353) - changed: E2, F2, G2, H2
Set-TextCell $wsQ4 "E2" "92.00"
Set-TextCell $wsQ4 "F2" "3.27"
Set-TextCell $wsQ4 "G2" "0.0111"
$wsQ4.Range("H2").Value = 8

# Row 3 (fund 008354) - changed: D3, E3, F3, G3, H3
Set-TextCell $wsQ4 "D3" "0.11"
Set-TextCell $wsQ4 "E3" "92.00"
Set-TextCell $wsQ4 "F3" "3.27"
Set-TextCell $wsQ4 "G3" "0.0036"
$wsQ4.Range("H3").Value = 8

# Row 4 (fund 002383) - brand new row
$wsQ4.Range("A4").Value = 2
Set-TextCell $wsQ4 "B4" "002383"
Set-TextCell $wsQ4 "C4" "大成趋势回报灵活配置混合"
Set-TextCell $wsQ4 "D4" "0.11"
Set-TextCell $wsQ4 "E4" "76.12"
Set-TextCell $wsQ4 "F4" "3.04"
Set-TextCell $wsQ4 "G4" "0.0033"
$wsQ4.Range("H4").Value = 6

# 3) Update the "总计" summary sheet: row 2 now reports 2022-Q4, and the old
#    2022-Q3 totals move down to a new row 3.
$oldCount = $wsTotal.Range("C2").Value2
$oldValue = $wsTotal.Range("D2").Value2

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = $oldCount
$wsTotal.Range("D3").Value = $oldValue

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.02

"done"
